$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item(1)

# C1: give it a top+bottom thin border (no left/right)
$ws1.Range("C1").ClearFormats()
$ws1.Range("C1").Borders.LineStyle = 1
$ws1.Range("C1").Borders(7).LineStyle = -4142
$ws1.Range("C1").Borders(10).LineStyle = -4142

# D1: give it a top+bottom+right thin border (no left)
$ws1.Range("D1").ClearFormats()
$ws1.Range("D1").Borders.LineStyle = 1
$ws1.Range("D1").Borders(7).LineStyle = -4142

# Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item(2)

# C1/D1 and F1/G1 need the exact same two border styles created above -
# copy them over (format only) instead of re-deriving each one, so the
# style table ends up with exactly the two new entries, no extras
$ws1.Range("C1").Copy()
$ws2.Range("C1").PasteSpecial(-4122)
$ws1.Range("D1").Copy()
$ws2.Range("D1").PasteSpecial(-4122)
$ws1.Range("C1").Copy()
$ws2.Range("F1").PasteSpecial(-4122)
$ws1.Range("D1").Copy()
$ws2.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Anonymize "fedcore" -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Drop the stray empty inline-string cell G5
$ws2.Range("G5").Clear()
